$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The row currently holding Patricia (004989791 / Patricia / 10941.72) is being
# split into two rows: Flavia and Daniela. Insert a new row above it so the
# original row's slot becomes Flavia's row and the newly inserted row becomes
# Daniela's row, pushing everything below down by one.
$ws.Rows.Item(3).Insert()

# Row 3 (was Patricia) -> Flavia
# Force the account-number column to be stored as text so the leading zeros
# in "004484207" are preserved instead of Excel auto-converting it to a number.
$ws.Cells.Item(3, 1).NumberFormat = "@"
$ws.Cells.Item(3, 1).Value = "004484207"
$ws.Cells.Item(3, 2).Value = "Flavia"
$ws.Cells.Item(3, 3).Value = 48900

# Row 4 (newly inserted blank row) -> Daniela
$ws.Cells.Item(4, 1).NumberFormat = "@"
$ws.Cells.Item(4, 1).Value = "004001621"
$ws.Cells.Item(4, 2).Value = "Daniela"
$ws.Cells.Item(4, 3).Value = 19331.42

# The Ana (005165116) row, previously row 5 with balance 7962.57, is now row 6
# after the insertion above. Update her balance to 9020.85.
$ws.Cells.Item(6, 3).Value = 9020.85
